$wb = $excel.ActiveWorkbook

# --- Sheet1: "singleCruiseTest" ---
$ws1 = $wb.Worksheets.Item("singleCruiseTest")
$ws1.Range("E6").Value = "N"
$ws1.Range("E6").Select()

# --- Sheet2: "parameterizedSearchTest" ---
$ws2 = $wb.Worksheets.Item("parameterizedSearchTest")

# Add new column G header
$ws2.Range("G1").Value = "runmode"
$ws2.Range("G1").HorizontalAlignment = -4108

# Fill in runmode values for existing rows
$ws2.Range("G2").Value = "N"
$ws2.Range("G2").HorizontalAlignment = -4108
$ws2.Range("G3").Value = "Y"
$ws2.Range("G3").HorizontalAlignment = -4108
$ws2.Range("G4").Value = "Y"
$ws2.Range("G4").HorizontalAlignment = -4108

# Fix F4 value (100000 -> 60000)
$ws2.Range("F4").Value = 60000

# Add new row 5
$ws2.Range("A5").Value = "BritishIslands"
$ws2.Range("B5").Value = "GreatBritain"
$ws2.Range("C5").Value = "Savona"
$ws2.Range("D5").Value = 7
$ws2.Range("E5").Value = "20 июнь"
$ws2.Range("F5").Value = 60000
$ws2.Range("G5").Value = "Y"
$ws2.Range("A5:G5").HorizontalAlignment = -4108

$ws2.Range("G2").Select()
